{"js": "// Make the word \"References\" bold while leaving the trailing period\n// un-bolded, splitting the single \"References.\" run into two runs.\nconst body = context.document.body;\n\nconst results = body.search(\"References\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('\"References\" not found in document body');\n}\n\n// Bold only the \"References\" text; Word automatically splits the run\n// so the un-searched \".\" keeps its original (non-bold) formatting.\nconst target = results.items[0];\ntarget.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Make the word \"References\" bold while leaving the trailing period\n# un-bolded. Word automatically splits the existing \"References.\" run\n# into a bold \"References\" run and a plain \".\" run when only part of\n# the run's formatting is changed.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"References\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n\n$found = $rng.Find.Execute()\nif (-not $found) {\n    throw '\"References\" not found in document content'\n}\n\n$rng.Font.Bold = 1\n"}
